# Apply the two textual corrections described in the commit.
$p = $ppt.ActivePresentation

# --- Slide 2: fix typo "pdf pf the commit" -> "pdf of the commit" -------
$slide2 = $p.Slides.Item(2)
$contentShape2 = $slide2.Shapes.Item(2)
$tr2 = $contentShape2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(2)
$run2c = $para2.Runs(3)
$run2c.Text = ". A pdf of the commit history is included in the zip file."

# --- Slide 3: merge the three runs of the final paragraph into one ------
$slide3 = $p.Slides.Item(3)
$contentShape3 = $slide3.Shapes.Item(2)
$tr3 = $contentShape3.TextFrame.TextRange
$para6 = $tr3.Paragraphs(6)

# Clearing a run's text removes it from the run sequence, shifting later
# runs down - so repeatedly clear the run that is now in position 2 to
# drop the old run2 then the old run3.
$para6.Runs(2).Text = ""
$para6b = $tr3.Paragraphs(6)
$para6b.Runs(2).Text = ""

# Finally, replace the remaining (first) run with the fully merged text.
$para6c = $tr3.Paragraphs(6)
$para6c.Runs(1).Text = "It then only remains to complete the remainder of the project in the most efficient and time effective manner possible."
